# Scheduled market-data refresh for Pandaemonium_Profits leve sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per row.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 348.16666
$ws.Range("I33").Value = 311.6154
$ws.Range("J33").Value = 443.2
$ws.Range("K33").Value = 311.6154
$ws.Range("L33").Value = 443.2
$ws.Range("M33").Value = -82.61540000000002
$ws.Range("N33").Value = -901.2
# Row 116
$ws.Range("H116").Value = 1986.9231
$ws.Range("I116").Value = 1748
$ws.Range("J116").Value = 2524.5
$ws.Range("K116").Value = 1748
$ws.Range("L116").Value = 2524.5
$ws.Range("M116").Value = 1694
$ws.Range("N116").Value = -9408.5
# Row 138
$ws.Range("H138").Value = 2649.3691
$ws.Range("I138").Value = 1097.3334
$ws.Range("J138").Value = 4440.1797
$ws.Range("K138").Value = 3292.0002
$ws.Range("L138").Value = 13320.5391
$ws.Range("M138").Value = 1847.9998
$ws.Range("N138").Value = -23600.5391

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1522.4
$ws.Range("I2").Value = 1804.3
$ws.Range("J2").Value = 958.6
$ws.Range("K2").Value = 1804.3
$ws.Range("L2").Value = 958.6
$ws.Range("M2").Value = -1691.3
$ws.Range("N2").Value = -1184.6
# Row 61
$ws.Range("H61").Value = 7281.48
$ws.Range("I61").Value = 3763.4707
$ws.Range("J61").Value = 14757.25
$ws.Range("K61").Value = 3763.4707
$ws.Range("L61").Value = 14757.25
$ws.Range("M61").Value = -3551.4707
$ws.Range("N61").Value = -15181.25
# Row 74
$ws.Range("H74").Value = 5739.393
$ws.Range("I74").Value = 2407.4736
$ws.Range("J74").Value = 12773.444
$ws.Range("K74").Value = 2407.4736
$ws.Range("L74").Value = 12773.444
$ws.Range("M74").Value = -1533.4736
$ws.Range("N74").Value = -14521.444
# Row 77
$ws.Range("H77").Value = 5739.393
$ws.Range("I77").Value = 2407.4736
$ws.Range("J77").Value = 12773.444
$ws.Range("K77").Value = 12037.368
$ws.Range("L77").Value = 63867.22
$ws.Range("M77").Value = -7669.367999999999
$ws.Range("N77").Value = -72603.22
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
# Row 116
$ws.Range("H116").Value = 1522.4
$ws.Range("I116").Value = 1804.3
$ws.Range("J116").Value = 958.6
$ws.Range("K116").Value = 1804.3
$ws.Range("L116").Value = 958.6
$ws.Range("M116").Value = 489.7
$ws.Range("N116").Value = -5546.6
# Row 136
$ws.Range("H136").Value = 7281.48
$ws.Range("I136").Value = 3763.4707
$ws.Range("J136").Value = 14757.25
$ws.Range("K136").Value = 11290.4121
$ws.Range("L136").Value = 44271.75
$ws.Range("M136").Value = -8740.4121
$ws.Range("N136").Value = -49371.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1522.4
$ws.Range("I3").Value = 1804.3
$ws.Range("J3").Value = 958.6
$ws.Range("K3").Value = 1804.3
$ws.Range("L3").Value = 958.6
$ws.Range("M3").Value = -1690.3
$ws.Range("N3").Value = -1186.6
# Row 64
$ws.Range("H64").Value = 396.66666
$ws.Range("I64").Value = 450
$ws.Range("J64").Value = 343.33334
$ws.Range("K64").Value = 450
$ws.Range("L64").Value = 343.33334
$ws.Range("M64").Value = -225
$ws.Range("N64").Value = -793.33334
# Row 67
$ws.Range("H67").Value = 396.66666
$ws.Range("I67").Value = 450
$ws.Range("J67").Value = 343.33334
$ws.Range("K67").Value = 450
$ws.Range("L67").Value = 343.33334
$ws.Range("M67").Value = 330
$ws.Range("N67").Value = -1903.33334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2090.224
$ws.Range("I31").Value = 1559.8948
$ws.Range("J31").Value = 3097.85
$ws.Range("K31").Value = 1559.8948
$ws.Range("L31").Value = 3097.85
$ws.Range("M31").Value = -1264.8948
$ws.Range("N31").Value = -3687.85
# Row 34
$ws.Range("H34").Value = 2090.224
$ws.Range("I34").Value = 1559.8948
$ws.Range("J34").Value = 3097.85
$ws.Range("K34").Value = 1559.8948
$ws.Range("L34").Value = 3097.85
$ws.Range("M34").Value = -1357.8948
$ws.Range("N34").Value = -3501.85
# Row 112
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
# Row 125
$ws.Range("H125").Value = 79790
$ws.Range("J125").Value = 79790
$ws.Range("L125").Value = 79790
$ws.Range("N125").Value = -84710
# Row 132
$ws.Range("H132").Value = 3498.1
$ws.Range("I132").Value = 3797.1904
$ws.Range("J132").Value = 2800.2222
$ws.Range("K132").Value = 11391.5712
$ws.Range("L132").Value = 8400.6666
$ws.Range("M132").Value = -8861.5712
$ws.Range("N132").Value = -13460.6666
# Row 134
$ws.Range("H134").Value = 2397.5078
$ws.Range("I134").Value = 1474.2162
$ws.Range("J134").Value = 3617.5715
$ws.Range("K134").Value = 4422.6486
$ws.Range("L134").Value = 10852.7145
$ws.Range("M134").Value = -1887.6486
$ws.Range("N134").Value = -15922.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 15662.308
$ws.Range("J131").Value = 17137.797
$ws.Range("L131").Value = 51413.391
$ws.Range("N131").Value = -61493.391

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3549.875
$ws.Range("I113").Value = 4933.3335
$ws.Range("J113").Value = 2719.8
$ws.Range("K113").Value = 4933.3335
$ws.Range("L113").Value = 2719.8
$ws.Range("M113").Value = -2763.3335
$ws.Range("N113").Value = -7059.8
# Row 132
$ws.Range("H132").Value = 4997.6
$ws.Range("I132").Value = 1972.8966
$ws.Range("J132").Value = 19617
$ws.Range("K132").Value = 5918.6898
$ws.Range("L132").Value = 58851
$ws.Range("M132").Value = -3388.6898
$ws.Range("N132").Value = -63911

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 8011960
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 10012450
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10012450
$ws.Range("M5").Value = -9888
$ws.Range("N5").Value = -10012674
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""
# Row 96
$ws.Range("H96").Value = 1002
$ws.Range("J96").Value = 1002
$ws.Range("L96").Value = 1002
$ws.Range("N96").Value = -3748
# Row 132
$ws.Range("H132").Value = 2474.1724
$ws.Range("I132").Value = 1179.8462
$ws.Range("J132").Value = 3525.8125
$ws.Range("K132").Value = 3539.5386
$ws.Range("L132").Value = 10577.4375
$ws.Range("M132").Value = -1009.5386
$ws.Range("N132").Value = -15637.4375
